$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 53-58 (A..J). Row 53 and 54 are overwritten in place,
# rows 55-58 are newly appended.
$rows = @(
    @{ r = 53; A = 2052; B = 1200; C = "Superior"; D = "Normal"; E = "A"; F = 1; G = 1; H = 1200; I = "Unoccupied"; J = "24-04-2020 03:28:54" },
    @{ r = 54; A = 2053; B = 1200; C = "Superior"; D = "Normal"; E = "A"; F = 1; G = 1; H = 1200; I = "Unoccupied"; J = "24-04-2020 03:30:52" },
    @{ r = 55; A = 2054; B = 1300; C = "Superior"; D = "Normal"; E = "A"; F = 1; G = 1; H = 1300; I = "Unoccupied"; J = "24-04-2020 03:31:01" },
    @{ r = 56; A = 2055; B = 1102; C = "Superior"; D = "Normal"; E = "A"; F = 1; G = 1; H = 1102; I = "Unoccupied"; J = "24-04-2020 03:33:59" },
    @{ r = 57; A = 2056; B = 9999; C = "Superior"; D = "Normal"; E = "A"; F = 1; G = 1; H = 9999; I = "Unoccupied"; J = "24-04-2020 03:34:24" },
    @{ r = 58; A = 2057; B = 9999; C = "Superior"; D = "Normal"; E = "A"; F = 1; G = 1; H = 9999; I = "Unoccupied"; J = "24-04-2020 03:35:03" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}
